$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" entries are plain decimal-looking text (e.g. "1.001") that must
# stay literal text exactly like the source feed, instead of being auto-parsed
# into a number by Excel's smart entry -- so format those cells as Text first.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.684.72"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.895.68"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "241.86"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4919"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").Value = "0.2938"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "0.06740"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").Value = "1.895.52"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "17.21"
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("D12").Value = "0.07240"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "90.77"
$ws.Range("E13").Value = "  +5.64%  "
$ws.Range("D14").Value = "0.6755"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "5.032"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "30.681.08"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "0.000007985"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "13.10"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").Value = "2.141.68"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "4.804"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "191.71"
$ws.Range("E23").Value = "  +33.75%  "
$ws.Range("D24").Value = "6.080"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").Value = "9.375"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Value = "156.14"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").Value = "19.04"
$ws.Range("E27").Value = "  +12.16%  "
$ws.Range("D28").Value = "1.898"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "1.404"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "4.298"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "0.09094"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").Value = "3.995"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "0.05202"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "0.7408"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("D35").Value = "1.107"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "2.761"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("D37").Value = "0.01832"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "0.9300"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "2.117"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "0.4393"
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("D42").Value = "105.14"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "5.735"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").Value = "7.532"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "0.1350"
$ws.Range("E46").Value = "  +5.35%  "
$ws.Range("D47").Value = "0.05863"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "8.757"
$ws.Range("E48").Value = "  +6.34%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.3926"
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "33.60"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.415"
$ws.Range("E51").Value = "  +6.03%  "
